# Update "想去人数" (wanted-attendance count) figures in the F column
# for the two sheets that carry this dataset: "展览" and "全部类型".
$wb = $excel.ActiveWorkbook

$updates = @{
    "F7"  = 37
    "F11" = 1459
    "F16" = 14
    "F23" = 3439
    "F24" = 415
    "F25" = 299
    "F26" = 471
    "F27" = 76
    "F28" = 23
    "F30" = 1181
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($cellRef in $updates.Keys) {
        $ws.Range($cellRef).Value = $updates[$cellRef]
    }
}
